$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new rows 11 and 12 by copying formatting from row 10 (gives correct A/B/C styles)
$ws.Range("A10:G10").Copy($ws.Range("A11:G11"))
$ws.Range("A10:G10").Copy($ws.Range("A12:G12"))

# Set all B-column (title) values first, in row order, so shared-string table fills title strings first
$ws.Range("B2").Value = 'Feldúlva találták a sírhelyeket'
$ws.Range("B3").Value = 'Elutasította medveügyben a Zetelaki Területtulajdonosi Társulás keresetét a táblabíróság'
$ws.Range("B4").Value = 'Visszatekintő: ezek voltak a közvéleményt leginkább foglalkoztató témáink 2019-ben'
$ws.Range("B5").Value = 'Ki akadályozta meg, hogy a barnamedve lekerüljön a szigorúan védett állatfajok listájáról Romániában?'
$ws.Range("B6").Value = 'Pénzt ígér a miniszter a medvék által veszélyeztetett települések védelmére'
$ws.Range("B7").Value = 'Hargita megye: megvan az év első 112-s medveészlelése'
$ws.Range("B8").Value = 'Aktívak a medvék Székelyudvarhely környékén'
$ws.Range("B9").Value = 'Lemondott a vadásztársaság az emberre támadó hidegkúti medve kilövéséről'
$ws.Range("B10").Value = 'Amíg elérhető közelségben van az ételmaradék, addig a medvék jelenlétére is számítani kell'
$ws.Range("B11").Value = 'Medveradar: Zetelakán és Farkaslakán voltak a legaktívabbak tavaly a nagyvadak'
$ws.Range("B12").Value = 'Gyergyószentmiklós utcáin kóborolt egy medve'

# Then set all C-column (content) values, in row order, so shared-string table appends content strings after
$ws.Range("C2").Value = 'Feldúlva találták a sírhelyeket, illetve azok környékét a Vasláb községhez tartozó hevederi temetőben. A nyomok alapján medvejárásra gyanakodnak.'
$ws.Range("C3").Value = 'Elutasította a Marosvásárhelyi Táblabíróság a Zetelaki Területtulajdonosi Társulás keresetét, amelyet a társulás a Környezetvédelmi Minisztérium ellen indított a vadgazdálkodási szabályozások alkalmazásának elmulasztása miatt. A társulás vezetője fellebbezést tervez.'
$ws.Range("C4").Value = 'Parkolás, pápalátogatás, temetőfoglalás, sportsikerek, véget nem érő medvetéma – a mögöttünk hagyott év székelyföldi közvéleményt leginkább érintő témáit gyűjtöttük csokorba, hónapokra lebontva.'
$ws.Range("C5").Value = 'Az állatvédők és az Európai Bizottság akadályozta meg, hogy a barnamedve öt évre lekerüljön a szigorúan védett állatfajok listájáról Romániában – állítja Benkő Erika RMDSZ-képviselő. '
$ws.Range("C6").Value = 'A háromszéki Zabolán tartott terepszemlét Costel Alexe környezetvédelmi miniszter, aki a látogatást követően arról számol be, hogy körvonalazódott egy olyan finanszírozási program, amely lehetővé teszi a medvék által veszélyeztetett települések védelmét.'
$ws.Range("C7").Value = 'Medvét látott a kertjében egy parajdi férfi a Sóhát utcában szombaton délután. A nagyvadat a gyümölcsfáknál fedezte fel, mintegy száz méterre a lakóháztól. Az esetet a 112-n jelentette, a helyszínre egy csendőri és egy mentőegység szállt ki.'
$ws.Range("C8").Value = 'Noha a magasabban fekvő térségekben már téli álmot alszanak a medvék, más területeken ez nem így van. Székelyudvarhely környékén például legalább tizenegy medve aktív jelenleg is, ezért a vadászok óvatosságra intenek.'
$ws.Range("C9").Value = 'Megúszta a kilövést a Hidegkúton emberre támadó medve, az illetékes vadásztársaságnál lemondtak arról, hogy a vad ártalmatlanítására rendkívüli jóváhagyást igényeljenek a környezetvédelmi minisztériumtól.'
$ws.Range("C10").Value = 'A szeméttárolók vonzzák a székelyudvarhelyi Cserehát lakónegyedbe az aktív nagyvadakat, ezért a Nagy-Küküllő Vadász- és Sporthorgász Egyesület medvebiztos kukákat rendelt, amelyeket a szemételszállító vállalattal egyeztetve helyezne ki. A medvék befogásával is próbálkoznak.'
$ws.Range("C11").Value = 'A kezdeti fellendülés után kissé lankadt az aktivitás a székelyföldi medveradar és -térkép néven emlegetett medveészlelő portálon, amelyet tavaly márciusban hozott létre Csala Dénes adatblogger. Ettől eltekintve a több mint 700 bejegyzést számláló medvetérkép a legszerteágazóbb adatbázisnak számít.'
$ws.Range("C12").Value = 'Egy városszéli üzemanyagtöltő állomásnál felbukkant medve miatt riasztották keddre virradóan a hatóságokat Gyergyószentmiklóson. Míg a csendőrök a nagyvadat kutatták, a vészhelyzeti felügyelőség a Ro-Alert rendszeren figyelmeztette a lakókat a veszélyre. A medvét végül megtalálták és elűzték.'

# Set remaining numeric columns (A id, D relevant, E severity, F deaths, G duplicate) and row heights
$ws.Range("A2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

$ws.Range("A3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Rows.Item(3).RowHeight = 105

$ws.Range("A4").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Rows.Item(4).RowHeight = 75

$ws.Range("A5").Value = 3
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("A6").Value = 4
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Rows.Item(6).RowHeight = 120

$ws.Range("A7").Value = 5
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Rows.Item(7).RowHeight = 90

$ws.Range("A8").Value = 6
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Rows.Item(8).RowHeight = 90

$ws.Range("A9").Value = 7
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Rows.Item(9).RowHeight = 90

$ws.Range("A10").Value = 8
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Rows.Item(10).RowHeight = 120

$ws.Range("A11").Value = 9
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Rows.Item(11).RowHeight = 120

$ws.Range("A12").Value = 10
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Rows.Item(12).RowHeight = 120

# Column widths for B and C (closest achievable quantized value to target 39.5703125)
$ws.Columns.Item(2).ColumnWidth = 38.65
$ws.Columns.Item(3).ColumnWidth = 38.65

# Update selection to match target view state (also clears topLeftCell scroll-anchor)
$ws.Range("E11").Select()

Write-Host "Done"